$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 292, shifting existing rows 292-331 down to 293-332.
$ws.Rows.Item(292).Insert()

# Populate the newly inserted row 292 with a new weekly price-record for Apio,
# matching the existing row's data pattern but with a new date (2023-07-17 -> serial 45124).
$ws.Range("A292").Value = 5
$ws.Range("B292").Value = "Macroferia Regional de Talca"
$ws.Range("C292").Value = "Maule"
$ws.Range("D292").Value = 45124
$ws.Range("E292").Value = 7
$ws.Range("F292").Value = 100112017
$ws.Range("G292").Value = "Apio"
$ws.Range("H292").Value = "Americana (o)"
$ws.Range("I292").Value = "Primera"
$ws.Range("J292").Value = 700
$ws.Range("K292").Value = 6000
$ws.Range("L292").Value = 6000
$ws.Range("M292").Value = 6000
$ws.Range("N292").Value = "`$/docena de matas"
$ws.Range("O292").Value = "Provincia del Elquí"
$ws.Range("P292").Value = 1000
$ws.Range("Q292").Value = 6
$ws.Range("R292").Value = "Hortaliza"
